# Generate Report for Handback
#
# Updates the localization-status report after a handback:
#  - Overview status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" (zh-cn & de-de columns).
#  - Per-language sheets (zh-cn, de-de) gain a "Latest Target File"
#    hyperlink (col I) and "Latest Handback File" (col J) value, and the
#    "Latest Handback DateTime" (col K) is stamped.
#  - Related columns are widened to fit the new, longer values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: Excel's ColumnWidth setter rounds through an integer pixel
# grid ( stored_width = round(ColumnWidth*6)/6 + 5/6 using the default
# Calibri-11 6px max-digit-width ). Pre-compensate so the persisted
# <col width="..."/> lands as close as possible to the real target.
# ---------------------------------------------------------------------
function Set-ColWidth($ws, $colIndex, $targetWidth) {
    $px = [Math]::Round(($targetWidth - 5.0/6.0) * 6.0)
    if ($px -lt 0) { $px = 0 }
    $input = $px / 6.0
    $ws.Columns.Item($colIndex).ColumnWidth = $input
}

# ---------------------------------------------------------------------
# 1. Overview sheet: handoff/handback status text + widen status cols
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

Set-ColWidth $overview 5 29.9777047293527
Set-ColWidth $overview 6 29.9777047293527

# ---------------------------------------------------------------------
# 2. Per-language detail sheets
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; TargetSuffix = "1dd3a6e17e2c1201dde9bc2f63a1d1f54fd84361.zh-cn.xlf"; HandbackDate = "2016-08-18 20:48:57" },
    @{ Name = "de-de"; TargetSuffix = "1dd3a6e17e2c1201dde9bc2f63a1d1f54fd84361.de-de.xlf"; HandbackDate = "2016-08-18 20:49:12" }
)

# Row 2 <-> 533a8485-... , Row 3 <-> 75af15d5-...
$rows = @(
    @{ Row = 2; Id = "533a8485-d1f1-499f-a28e-12dbe2b027d5"; TargetHash = "1dd3a6e17e2c1201dde9bc2f63a1d1f54fd84361" },
    @{ Row = 3; Id = "75af15d5-9871-42fc-9627-421b87f9cd98"; TargetHash = "b9172a4d956fa921458523b96750f6c60140bc4d" }
)

foreach ($langInfo in $langSheets) {
    $ws = $wb.Worksheets.Item($langInfo.Name)

    # "Status" (C) moves from "Ready for handoff" to the handback state,
    # same text as the Overview sheet (both share the same underlying
    # shared string in the source workbook).
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Widen "Status" (C) and the new "Latest Target File"/"Latest Handback
    # File" (I/J) columns so the longer values/hyperlinks fit.
    Set-ColWidth $ws 3 29.9777047293527
    Set-ColWidth $ws 9 40
    Set-ColWidth $ws 10 40

    foreach ($rowInfo in $rows) {
        $r = $rowInfo.Row
        $mdFile = $rowInfo.Id + ".md"
        $xlfFile = $rowInfo.Id + "." + $rowInfo.TargetHash + "." + $langInfo.Name + ".xlf"
        $mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb3b58a3e98790a53ba9044c088470a3af38b2fa/e2e/" + $mdFile

        # Column I: "Latest Target File" -> hyperlink to the same source doc
        $cellI = $ws.Range("I" + $r)
        $ws.Hyperlinks.Add($cellI, $mdUrl, "", "", $mdFile) | Out-Null
        $cellI.Font.Underline = $true
        $cellI.Font.Color = 15570276

        # Column J: "Latest Handback File" -> generated xliff file name
        $ws.Range("J" + $r).Value = $xlfFile

        # Column K: "Latest Handback DateTime"
        $ws.Range("K" + $r).Value = $langInfo.HandbackDate
    }
}
